# Update result values pulled from server for each year sheet.
$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 973.9537847600009
$ws2025.Range("E2").Value = 28982.37596598056
$ws2025.Range("I2").Value = 16175.28135478
$ws2025.Range("L2").Value = 48524.529503538
$ws2025.Range("M2").Value = 10590.587968015
$ws2025.Range("N2").Value = 7166.934239853319
$ws2025.Range("O2").Value = 6983.506508100742

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 5712.560177842886
$ws2030.Range("E2").Value = 56106.05588781912
$ws2030.Range("I2").Value = 44217.8984721661
$ws2030.Range("L2").Value = 66966.57749858923
$ws2030.Range("M2").Value = 21984.28023276101
$ws2030.Range("N2").Value = 10610.95466116662
$ws2030.Range("O2").Value = 12070.09906536593

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 2861.961401238371
$ws2035.Range("B2").Value = 8026.889663087295
$ws2035.Range("E2").Value = 67297.73995507321
$ws2035.Range("I2").Value = 59256.42575923612
$ws2035.Range("L2").Value = 66966.57749858923
$ws2035.Range("M2").Value = 25464.6214365565
$ws2035.Range("N2").Value = 15155.33304031695
$ws2035.Range("O2").Value = 14768.69820138149

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 2861.961401238371
$ws2040.Range("B2").Value = 8026.889663087295
$ws2040.Range("E2").Value = 67297.73995507321
$ws2040.Range("I2").Value = 59256.42575923612
$ws2040.Range("L2").Value = 66966.57749858923
$ws2040.Range("M2").Value = 25464.6214365565
$ws2040.Range("N2").Value = 15260.17788384349
$ws2040.Range("O2").Value = 14768.69820138149

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 6302.873118834019
$ws2045.Range("B2").Value = 8026.889663087295
$ws2045.Range("E2").Value = 67297.73995507321
$ws2045.Range("I2").Value = 59256.42575923612
$ws2045.Range("L2").Value = 66966.57749858923
$ws2045.Range("M2").Value = 25464.6214365565
$ws2045.Range("N2").Value = 15798.2207999018
$ws2045.Range("O2").Value = 17109.77519465316

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Value = 6302.873118834019
$ws2050.Range("B2").Value = 8026.889663087295
$ws2050.Range("E2").Value = 67297.73995507321
$ws2050.Range("I2").Value = 59256.42575923612
$ws2050.Range("L2").Value = 66966.57749858923
$ws2050.Range("M2").Value = 25464.6214365565
$ws2050.Range("N2").Value = 15798.2207999018
$ws2050.Range("O2").Value = 17109.77519465316
